$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header text "MODEL_CONDITION" becomes "MODELCONDITION"
[void]$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION")

# Column A (the per-row index numbers 1/3/10/18 duplicated later as the
# last column) is removed entirely; every remaining column shifts one
# place to the left (old B->A, C->B, D->C, E->D, F->E), so the sheet
# ends up spanning A1:E5 instead of A1:F5.
$ws.Columns.Item(1).Delete()
